$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 02:03"

# Update country stat rows that changed (new data pulled in; some countries
# overtook others in the ranking so their rows swapped places)

# Row 4
$ws.Cells.Item(4, 2).Value = 1262692
$ws.Cells.Item(4, 3).Value = 25059
$ws.Cells.Item(4, 4).Value = 206308
$ws.Cells.Item(4, 5).Value = 981593
$ws.Cells.Item(4, 7).Value = 2520
$ws.Cells.Item(4, 8).Value = 74791

# Row 12
$ws.Cells.Item(12, 2).Value = 126148
$ws.Cells.Item(12, 3).Value = 11433
$ws.Cells.Item(12, 4).Value = 51370
$ws.Cells.Item(12, 5).Value = 66212
$ws.Cells.Item(12, 7).Value = 645
$ws.Cells.Item(12, 8).Value = 8566

# Row 15
$ws.Cells.Item(15, 2).Value = 63496
$ws.Cells.Item(15, 3).Value = 1450
$ws.Cells.Item(15, 5).Value = 31448
$ws.Cells.Item(15, 7).Value = 189
$ws.Cells.Item(15, 8).Value = 4232

# Row 20
$ws.Cells.Item(20, 6).Value = 137

# Row 46
$ws.Cells.Item(46, 4).Value = 1960
$ws.Cells.Item(46, 5).Value = 6485

# Row 57
$ws.Cells.Item(57, 2).Value = 5208
$ws.Cells.Item(57, 3).Value = 188
$ws.Cells.Item(57, 5).Value = 3411
$ws.Cells.Item(57, 7).Value = 9
$ws.Cells.Item(57, 8).Value = 273

# Row 64
$ws.Cells.Item(64, 1).Value = "Nigeria"
$ws.Cells.Item(64, 2).Value = 3145
$ws.Cells.Item(64, 3).Value = 195
$ws.Cells.Item(64, 4).Value = 534
$ws.Cells.Item(64, 5).Value = 2508
$ws.Cells.Item(64, 6).Value = 4
$ws.Cells.Item(64, 7).Value = 5
$ws.Cells.Item(64, 8).Value = 103

# Row 65
$ws.Cells.Item(65, 1).Value = "Hungria"
$ws.Cells.Item(65, 2).Value = 3111
$ws.Cells.Item(65, 3).Value = 46
$ws.Cells.Item(65, 4).Value = 759
$ws.Cells.Item(65, 5).Value = 1979
$ws.Cells.Item(65, 6).Value = 50
$ws.Cells.Item(65, 7).Value = 10
$ws.Cells.Item(65, 8).Value = 373

# Row 66
$ws.Cells.Item(66, 1).Value = "Tailandia"
$ws.Cells.Item(66, 2).Value = 2989
$ws.Cells.Item(66, 3).Value = 1
$ws.Cells.Item(66, 4).Value = 2761
$ws.Cells.Item(66, 5).Value = 173
$ws.Cells.Item(66, 6).Value = 61
$ws.Cells.Item(66, 7).Value = 1
$ws.Cells.Item(66, 8).Value = 55

# Row 72
$ws.Cells.Item(72, 5).Value = 1157
$ws.Cells.Item(72, 7).Value = 44
$ws.Cells.Item(72, 8).Value = 108

# Row 82
$ws.Cells.Item(82, 6).Value = 14

# Row 93
$ws.Cells.Item(93, 2).Value = 1025
$ws.Cells.Item(93, 3).Value = 3
$ws.Cells.Item(93, 4).Value = 591
$ws.Cells.Item(93, 5).Value = 391
$ws.Cells.Item(93, 6).Value = 18

# Row 110
$ws.Cells.Item(110, 2).Value = 673
$ws.Cells.Item(110, 3).Value = 3
$ws.Cells.Item(110, 4).Value = 486
$ws.Cells.Item(110, 5).Value = 170

# Row 154
$ws.Cells.Item(154, 1).Value = "Bermudas"
$ws.Cells.Item(154, 2).Value = 118
$ws.Cells.Item(154, 3).Value = 3
$ws.Cells.Item(154, 4).Value = 59
$ws.Cells.Item(154, 5).Value = 52
$ws.Cells.Item(154, 6).Value = 4
$ws.Cells.Item(154, 8).Value = 7

# Row 155
$ws.Cells.Item(155, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(155, 2).Value = 116
$ws.Cells.Item(155, 4).Value = 103
$ws.Cells.Item(155, 5).Value = 5
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 8).Value = 8

# Row 164
$ws.Cells.Item(164, 2).Value = 92
$ws.Cells.Item(164, 3).Value = 3
$ws.Cells.Item(164, 5).Value = 55

# Row 167
$ws.Cells.Item(167, 4).Value = 21
$ws.Cells.Item(167, 5).Value = 60

# Row 209
$ws.Cells.Item(209, 1).Value = "Comoras"
$ws.Cells.Item(209, 3).Value = 5
$ws.Cells.Item(209, 4).Value = 0
$ws.Cells.Item(209, 5).Value = 7
$ws.Cells.Item(209, 7).Value = 1

# Row 210
$ws.Cells.Item(210, 1).Value = "Mauritania"
$ws.Cells.Item(210, 4).Value = 6
$ws.Cells.Item(210, 5).Value = 1
$ws.Cells.Item(210, 8).Value = 1

# Row 211
$ws.Cells.Item(211, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(211, 2).Value = 8
$ws.Cells.Item(211, 4).Value = 8
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 8).Value = 0

# Row 212
$ws.Cells.Item(212, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(212, 4).Value = 3
$ws.Cells.Item(212, 5).Value = 3
$ws.Cells.Item(212, 8).Value = 1

# Row 213
$ws.Cells.Item(213, 1).Value = "Butan"
$ws.Cells.Item(213, 2).Value = 7
$ws.Cells.Item(213, 4).Value = 5
$ws.Cells.Item(213, 5).Value = 2

# Row 214
$ws.Cells.Item(214, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(214, 4).Value = 0
$ws.Cells.Item(214, 5).Value = 6

# Row 215
$ws.Cells.Item(215, 1).Value = "Sahara Occidental"
$ws.Cells.Item(215, 4).Value = 5
$ws.Cells.Item(215, 5).Value = 1

# Row 216
$ws.Cells.Item(216, 1).Value = "San Bartolome"
$ws.Cells.Item(216, 2).Value = 6
$ws.Cells.Item(216, 3).Value = 0
$ws.Cells.Item(216, 4).Value = 6
$ws.Cells.Item(216, 5).Value = 0
$ws.Cells.Item(216, 7).Value = 0
$ws.Cells.Item(216, 8).Value = 0
